# Update master_solution_decisions sheet with the Mixed Logit Model results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Aarau -> Chiasso, Road)
$ws.Range("E2").Value = 6.342146727268994
$ws.Range("G2").Value = 66.15000000000001

# Row 3 (Aarau -> Chiasso, Intermodal)
$ws.Range("E3").Value = 1.010625
$ws.Range("G3").Value = 80.84999999999999

# Row 4 (Chiasso -> Aarau, Road)
$ws.Range("F4").Value = 392.0444

# Row 5 (Chiasso -> Aarau, Intermodal)
$ws.Range("F5").Value = 248.762442

# Rows 6-9: x_open flips from 0 to -0 (sign-only change on a zero value)
$ws.Range("D6").Value = -0
$ws.Range("D7").Value = -0
$ws.Range("D8").Value = -0
$ws.Range("D9").Value = -0

# Row 14 (Visp -> Chiasso, Road)
$ws.Range("E14").Value = 18.68936771379288
$ws.Range("G14").Value = 1693.440000000003

# Row 15 (Visp -> Chiasso, Intermodal)
$ws.Range("G15").Value = 70.5600000000007

# Row 16 (Chiasso -> Visp, Road)
$ws.Range("E16").Value = 18.61888993138674
$ws.Range("G16").Value = 1384.320000000001

# Row 17 (Chiasso -> Visp, Intermodal)
$ws.Range("G17").Value = 57.68
